# "Generate Report for Archive" — the localization status report was
# regenerated: the zh-cn / de-de handoff status moved on from
# "Ready for handoff" to "In Translation". That status string is echoed on
# the "Overview" sheet (columns E/F, one per language) as well as on each
# language's own sheet ("zh-cn" / "de-de", column C). Re-generating the
# report also re-ran the report tool's column autosize pass, which narrowed
# the Status column now that the new text is shorter than the old text.

$wb = $excel.ActiveWorkbook

# Overview sheet: the zh-cn (E) and de-de (F) status cells on row 2.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsOverview.Range("E1:E2").ColumnWidth = 12.5
$wsOverview.Range("F1:F2").ColumnWidth = 12.5

# zh-cn sheet: Status column (C)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"
$wsZhCn.Range("C1:C2").ColumnWidth = 12.5

# de-de sheet: Status column (C)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"
$wsDeDe.Range("C1:C2").ColumnWidth = 12.5
